# Add two new result sheets (VET_persistence, VET_extrapolation) at the end
# of the workbook, matching the "VET" (variational echo tracking) nowcasting
# method results that were appended alongside the existing LK / DARTS /
# proesmans sheets.

$wb = $excel.ActiveWorkbook

# Preserve the author's last interactive selection on proesmans_anvil
# (cell F23) before the new sheets take over as the active tab.
$lastActive = $wb.Worksheets.Item("proesmans_anvil")
$lastActive.Activate()
$lastActive.Range("F23").Select()

# Use an existing results sheet as the template for the bold/centered/
# bordered header-row formatting used throughout the workbook.
$headerTemplate = $wb.Worksheets.Item("LK_persistence").Range("A1:G1")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$ws1 = $wb.Worksheets.Add($null, $lastSheet)
$ws1.Name = "VET_persistence"
$headerTemplate.Copy()
$ws1.Range("A1:G1").PasteSpecial(-4122)
$ws1.PageSetup.LeftMargin = 54
$ws1.PageSetup.RightMargin = 54
$ws1.PageSetup.TopMargin = 72
$ws1.PageSetup.BottomMargin = 72
$ws1.PageSetup.HeaderMargin = 36
$ws1.PageSetup.FooterMargin = 36

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "VET_extrapolation"
$headerTemplate.Copy()
$ws2.Range("A1:G1").PasteSpecial(-4122)
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36


# --- VET_persistence data ---
$ws1.Range("A1").Value = "Event"
$ws1.Range("B1").Value = "Lead Time 1"
$ws1.Range("C1").Value = "Lead Time 2"
$ws1.Range("D1").Value = "Lead Time 3"
$ws1.Range("E1").Value = "Lead Time 4"
$ws1.Range("F1").Value = "Lead Time 5"
$ws1.Range("G1").Value = "Lead Time 6"
$ws1.Range("A2").Value = "20230520_2235"
$ws1.Range("B2").Value = 0.249247133292094
$ws1.Range("C2").Value = 0.12906601699843
$ws1.Range("D2").Value = 0.08269537310337292
$ws1.Range("E2").Value = 0.05488475457633449
$ws1.Range("F2").Value = 0.03967948456035416
$ws1.Range("G2").Value = 0.03059399439466735
$ws1.Range("A3").Value = "20190320_0005"
$ws1.Range("B3").Value = 0.2782406947954387
$ws1.Range("C3").Value = 0.1409440101977913
$ws1.Range("D3").Value = 0.08221309531284574
$ws1.Range("E3").Value = 0.05193464064243412
$ws1.Range("F3").Value = 0.03331682066020348
$ws1.Range("G3").Value = 0.01694535906508861
$ws1.Range("A4").Value = "20191222_0900"
$ws1.Range("B4").Value = 0.04365638963673633
$ws1.Range("C4").Value = 0.021352653638231
$ws1.Range("D4").Value = 0.01704859472317032
$ws1.Range("E4").Value = 0.01482767664884353
$ws1.Range("F4").Value = 0.02487907511365598
$ws1.Range("G4").Value = 0.03124569686151027
$ws1.Range("A5").Value = "20180505_1745"
$ws1.Range("B5").Value = 0.3665606431991816
$ws1.Range("C5").Value = 0.2246164583985399
$ws1.Range("D5").Value = 0.1570187829568649
$ws1.Range("E5").Value = 0.1088552528422098
$ws1.Range("F5").Value = 0.07376949583319505
$ws1.Range("G5").Value = 0.04974047026649884
$ws1.Range("A6").Value = "20230513_1455"
$ws1.Range("B6").Value = 0.04156588623373377
$ws1.Range("C6").Value = 0.04740432508664978
$ws1.Range("D6").Value = 0.02690902791477917
$ws1.Range("E6").Value = 0.0314833300476104
$ws1.Range("F6").Value = 0.0324439671008708
$ws1.Range("G6").Value = 0.01253703578818633
$ws1.Range("A7").Value = "20200911_1315"
$ws1.Range("B7").Value = 0.4714376744476271
$ws1.Range("C7").Value = 0.2804363077581025
$ws1.Range("D7").Value = 0.1846323716590821
$ws1.Range("E7").Value = 0.1285092895410339
$ws1.Range("F7").Value = 0.08763271036740511
$ws1.Range("G7").Value = 0.05808281344603691
$ws1.Range("A8").Value = "20191111_0710"
$ws1.Range("B8").Value = 0.1111474737964534
$ws1.Range("C8").Value = 0.0592039448958813
$ws1.Range("D8").Value = 0.04320225739104044
$ws1.Range("E8").Value = 0.03448699904862525
$ws1.Range("F8").Value = 0.0369132860298621
$ws1.Range("G8").Value = 0.02369487936731393
$ws1.Range("A9").Value = "20230302_0245"
$ws1.Range("B9").Value = 0.08468256043233927
$ws1.Range("C9").Value = 0.02111070794461422
$ws1.Range("D9").Value = 0.04028291821984258
$ws1.Range("E9").Value = 0.04741679070896806
$ws1.Range("F9").Value = 0.04510873280297124
$ws1.Range("G9").Value = 0.04313551231357516
$ws1.Range("A10").Value = "20190412_1220"
$ws1.Range("B10").Value = 0.2935024979585278
$ws1.Range("C10").Value = 0.1441972377725796
$ws1.Range("D10").Value = 0.0745058114340333
$ws1.Range("E10").Value = 0.03913489122854651
$ws1.Range("F10").Value = 0.02479409266691508
$ws1.Range("G10").Value = 0.01832741397874997
$ws1.Range("A11").Value = "20200120_1440"
$ws1.Range("B11").Value = 0.02251070938578288
$ws1.Range("C11").Value = 0.00913848103339623
$ws1.Range("D11").Value = 0.01071753188920646
$ws1.Range("E11").Value = 0.008470671907106087
$ws1.Range("F11").Value = 0.005416906669745094
$ws1.Range("G11").Value = 0.01480975717137301
$ws1.Range("A12").Value = "20230129_2215"
$ws1.Range("B12").Value = 0.10208822967322
$ws1.Range("C12").Value = 0.04567519337518174
$ws1.Range("D12").Value = 0.02224271994897925
$ws1.Range("E12").Value = 0.009123456545722082
$ws1.Range("F12").Value = 0.007030939446932111
$ws1.Range("G12").Value = 0.006970795233397075
$ws1.Range("A13").Value = "20181014_0515"
$ws1.Range("B13").Value = 0.4458698825363287
$ws1.Range("C13").Value = 0.3016227598989668
$ws1.Range("D13").Value = 0.2324808076399397
$ws1.Range("E13").Value = 0.1928626935923116
$ws1.Range("F13").Value = 0.1699279894102941
$ws1.Range("G13").Value = 0.1540260749425108

# --- VET_extrapolation data ---
$ws2.Range("A1").Value = "Event"
$ws2.Range("B1").Value = "Lead Time 1"
$ws2.Range("C1").Value = "Lead Time 2"
$ws2.Range("D1").Value = "Lead Time 3"
$ws2.Range("E1").Value = "Lead Time 4"
$ws2.Range("F1").Value = "Lead Time 5"
$ws2.Range("G1").Value = "Lead Time 6"
$ws2.Range("A2").Value = "20230520_2235"
$ws2.Range("B2").Value = 0.420651081612502
$ws2.Range("C2").Value = 0.2725751641508029
$ws2.Range("D2").Value = 0.1956241067337803
$ws2.Range("E2").Value = 0.1496992567654788
$ws2.Range("F2").Value = 0.1227726870014735
$ws2.Range("G2").Value = 0.1041607561076801
$ws2.Range("A3").Value = "20190320_0005"
$ws2.Range("B3").Value = 0.5154929260055031
$ws2.Range("C3").Value = 0.3597764085069424
$ws2.Range("D3").Value = 0.2701540848191289
$ws2.Range("E3").Value = 0.2077461944812471
$ws2.Range("F3").Value = 0.1569357652122197
$ws2.Range("G3").Value = 0.1148136699116411
$ws2.Range("A4").Value = "20191222_0900"
$ws2.Range("B4").Value = 0.02626277391822984
$ws2.Range("C4").Value = 0.009817834749130071
$ws2.Range("D4").Value = 0.008782901149979082
$ws2.Range("E4").Value = 0.005513699153356793
$ws2.Range("F4").Value = 0.009433258124923946
$ws2.Range("G4").Value = 0.001762965118971555
$ws2.Range("A5").Value = "20180505_1745"
$ws2.Range("B5").Value = 0.5099904241144766
$ws2.Range("C5").Value = 0.3301444459352756
$ws2.Range("D5").Value = 0.2334948942140488
$ws2.Range("E5").Value = 0.174222289633173
$ws2.Range("F5").Value = 0.1356305237470961
$ws2.Range("G5").Value = 0.1104820346305269
$ws2.Range("A6").Value = "20230513_1455"
$ws2.Range("B6").Value = 0.1404049397650292
$ws2.Range("C6").Value = 0.05629383311269984
$ws2.Range("D6").Value = 0.02137234829575282
$ws2.Range("E6").Value = 0.009037475357010534
$ws2.Range("F6").Value = 0.002269355255026133
$ws2.Range("G6").Value = 0.00004972289050637027
$ws2.Range("A7").Value = "20200911_1315"
$ws2.Range("B7").Value = 0.5447077939689181
$ws2.Range("C7").Value = 0.3465703514516474
$ws2.Range("D7").Value = 0.2406294579771947
$ws2.Range("E7").Value = 0.1781208174968559
$ws2.Range("F7").Value = 0.1366315329948568
$ws2.Range("G7").Value = 0.1071345243756502
$ws2.Range("A8").Value = "20191111_0710"
$ws2.Range("B8").Value = 0.2975649768502179
$ws2.Range("C8").Value = 0.1458342134935508
$ws2.Range("D8").Value = 0.07362327641002393
$ws2.Range("E8").Value = 0.04367958292779116
$ws2.Range("F8").Value = 0.02995042094263866
$ws2.Range("G8").Value = 0.0221138918323158
$ws2.Range("A9").Value = "20230302_0245"
$ws2.Range("B9").Value = 0.2329184717031176
$ws2.Range("C9").Value = 0.1026447133605268
$ws2.Range("D9").Value = 0.04946216455639678
$ws2.Range("E9").Value = 0.0282526825854036
$ws2.Range("F9").Value = 0.01401882788713467
$ws2.Range("G9").Value = 0.004743102772867889
$ws2.Range("A10").Value = "20190412_1220"
$ws2.Range("B10").Value = 0.449891428162017
$ws2.Range("C10").Value = 0.259844830738637
$ws2.Range("D10").Value = 0.1635656847750692
$ws2.Range("E10").Value = 0.1038713131841534
$ws2.Range("F10").Value = 0.07559447595469238
$ws2.Range("G10").Value = 0.06363471993739828
$ws2.Range("A11").Value = "20200120_1440"
$ws2.Range("B11").Value = 0.1036433349586603
$ws2.Range("C11").Value = 0.04615471236984033
$ws2.Range("D11").Value = 0.02127763683914593
$ws2.Range("E11").Value = 0.009344732407291102
$ws2.Range("F11").Value = 0.004215479953717881
$ws2.Range("G11").Value = 0.006159698073101816
$ws2.Range("A12").Value = "20230129_2215"
$ws2.Range("B12").Value = 0.1632286531960156
$ws2.Range("C12").Value = 0.04737332421031459
$ws2.Range("D12").Value = 0.01726621367457788
$ws2.Range("E12").Value = 0.008966109408683315
$ws2.Range("F12").Value = 0.004347749225487287
$ws2.Range("G12").Value = 0.001983759454499204
$ws2.Range("A13").Value = "20181014_0515"
$ws2.Range("B13").Value = 0.5360559231426352
$ws2.Range("C13").Value = 0.3880676569803455
$ws2.Range("D13").Value = 0.303706360109726
$ws2.Range("E13").Value = 0.2487669411507497
$ws2.Range("F13").Value = 0.2098501251083696
$ws2.Range("G13").Value = 0.1794197813127008
